$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 17966.666
$ws.Range("I70").Value = 20960
$ws.Range("K70").Value = 62880
$ws.Range("M70").Value = -62610
$ws.Range("H73").Value = 17966.666
$ws.Range("I73").Value = 20960
$ws.Range("K73").Value = 62880
$ws.Range("M73").Value = -61944
$ws.Range("H111").Value = 2375.0557
$ws.Range("I111").Value = 950.4545000000001
$ws.Range("J111").Value = 4613.7144
$ws.Range("K111").Value = 2851.3635
$ws.Range("L111").Value = 13841.1432
$ws.Range("M111").Value = 215.6364999999996
$ws.Range("N111").Value = -19975.1432
$ws.Range("H137").Value = 40160.42
$ws.Range("I137").Value = 1299.0769
$ws.Range("J137").Value = 79021.766
$ws.Range("K137").Value = 3897.2307
$ws.Range("L137").Value = 237065.298
$ws.Range("M137").Value = -1347.2307
$ws.Range("N137").Value = -242165.298
$ws.Range("H138").Value = 2495.35
$ws.Range("I138").Value = 633.375
$ws.Range("J138").Value = 2960.8438
$ws.Range("K138").Value = 1900.125
$ws.Range("L138").Value = 8882.5314
$ws.Range("M138").Value = 3239.875
$ws.Range("N138").Value = -19162.5314

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 25478.4
$ws.Range("I32").Value = 30629
$ws.Range("K32").Value = 30629
$ws.Range("M32").Value = -30342
$ws.Range("H43").Value = 24725.4
$ws.Range("J43").Value = 20959
$ws.Range("L43").Value = 20959
$ws.Range("N43").Value = -21585
$ws.Range("H45").Value = 2729.6191
$ws.Range("I45").Value = 2385.875
$ws.Range("K45").Value = 2385.875
$ws.Range("M45").Value = -2008.875
$ws.Range("H88").Value = 42580.16
$ws.Range("I88").Value = 1394.5
$ws.Range("J88").Value = 61961.65
$ws.Range("K88").Value = 1394.5
$ws.Range("L88").Value = 61961.65
$ws.Range("M88").Value = -988.5
$ws.Range("N88").Value = -62773.65
$ws.Range("H91").Value = 42580.16
$ws.Range("I91").Value = 1394.5
$ws.Range("J91").Value = 61961.65
$ws.Range("K91").Value = 1394.5
$ws.Range("L91").Value = 61961.65
$ws.Range("M91").Value = 9.5
$ws.Range("N91").Value = -64769.65
$ws.Range("H101").Value = 50000
$ws.Range("J101").Value = 50000
$ws.Range("L101").Value = 50000
$ws.Range("N101").Value = -56490
$ws.Range("H102").Value = 2959.6191
$ws.Range("I102").Value = 2267.4285
$ws.Range("J102").Value = 3305.7144
$ws.Range("K102").Value = 2267.4285
$ws.Range("L102").Value = 3305.7144
$ws.Range("M102").Value = -645.4285
$ws.Range("N102").Value = -6549.7144
$ws.Range("H119").Value = 50000
$ws.Range("J119").Value = 50000
$ws.Range("L119").Value = 50000
$ws.Range("N119").Value = -59676
$ws.Range("H122").Value = 2249.9524
$ws.Range("I122").Value = 1997.6666
$ws.Range("J122").Value = 2586.3333
$ws.Range("K122").Value = 5992.9998
$ws.Range("L122").Value = 7758.999899999999
$ws.Range("M122").Value = -3542.9998
$ws.Range("N122").Value = -12658.9999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 32899.938
$ws.Range("J86").Value = 2599.8
$ws.Range("L86").Value = 2599.8
$ws.Range("N86").Value = -4845.8
$ws.Range("H89").Value = 32899.938
$ws.Range("J89").Value = 2599.8
$ws.Range("L89").Value = 12999
$ws.Range("N89").Value = -24231
$ws.Range("H105").Value = 1973.3334
$ws.Range("J105").Value = 2088.889
$ws.Range("L105").Value = 2088.889
$ws.Range("N105").Value = -5582.889

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4802.8667
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").ClearContents()
$ws.Range("H34").Value = 4802.8667
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").ClearContents()
$ws.Range("H41").Value = 20725
$ws.Range("J41").Value = 26450
$ws.Range("L41").Value = 26450
$ws.Range("N41").Value = -27306
$ws.Range("H59").Value = 21000
$ws.Range("J59").Value = 21000
$ws.Range("L59").Value = 21000
$ws.Range("N59").Value = -23290
$ws.Range("H99").Value = 5236.143
$ws.Range("I99").Value = 3766.077
$ws.Range("K99").Value = 3766.077
$ws.Range("M99").Value = -2268.077
$ws.Range("H126").Value = 5236.143
$ws.Range("I126").Value = 3766.077
$ws.Range("K126").Value = 11298.231
$ws.Range("M126").Value = -8828.231

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H58").Value = 3276.75
$ws.Range("J58").Value = 4004.8
$ws.Range("L58").Value = 12014.4
$ws.Range("N58").Value = -12270.4
$ws.Range("H68").Value = 1287.7084
$ws.Range("J68").Value = 1345.7273
$ws.Range("L68").Value = 4037.1819
$ws.Range("N68").Value = -5659.1819
$ws.Range("H71").Value = 1287.7084
$ws.Range("J71").Value = 1345.7273
$ws.Range("L71").Value = 12111.5457
$ws.Range("N71").Value = -20223.5457
$ws.Range("H86").Value = 805.2
$ws.Range("J86").Value = 725
$ws.Range("L86").Value = 2175
$ws.Range("N86").Value = -4547
$ws.Range("H89").Value = 805.2
$ws.Range("J89").Value = 725
$ws.Range("L89").Value = 6525
$ws.Range("N89").Value = -18381
$ws.Range("H107").Value = 4602.68
$ws.Range("I107").Value = 9574.091
$ws.Range("J107").Value = 696.5714
$ws.Range("K107").Value = 28722.273
$ws.Range("L107").Value = 2089.7142
$ws.Range("M107").Value = -26802.273
$ws.Range("N107").Value = -5929.7142
$ws.Range("H113").Value = 394.66666
$ws.Range("I113").Value = 380
$ws.Range("J113").Value = 402
$ws.Range("K113").Value = 1140
$ws.Range("L113").Value = 1206
$ws.Range("M113").Value = 1030
$ws.Range("N113").Value = -5546
$ws.Range("H131").Value = 119890.82
$ws.Range("J131").Value = 127416.914
$ws.Range("L131").Value = 382250.742
$ws.Range("N131").Value = -392330.742

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 81173.42
$ws.Range("I132").Value = 78499.766
$ws.Range("J132").Value = 86966.336
$ws.Range("K132").Value = 235499.298
$ws.Range("L132").Value = 260899.008
$ws.Range("M132").Value = -232969.298
$ws.Range("N132").Value = -265959.008

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1133.3636
$ws.Range("I16").Value = 718.55554
$ws.Range("J16").Value = 3000
$ws.Range("K16").Value = 718.55554
$ws.Range("L16").Value = 3000
$ws.Range("M16").Value = -548.55554
$ws.Range("N16").Value = -3340
$ws.Range("H22").Value = 1500
$ws.Range("I22").Value = 500
$ws.Range("K22").Value = 500
$ws.Range("M22").Value = -205
$ws.Range("H27").Value = 1500
$ws.Range("I27").Value = 500
$ws.Range("K27").Value = 500
$ws.Range("M27").Value = -393
$ws.Range("H82").Value = 2219.5
$ws.Range("I82").Value = 2355
$ws.Range("K82").Value = 2355
$ws.Range("M82").Value = -1994
$ws.Range("H85").Value = 2219.5
$ws.Range("I85").Value = 2355
$ws.Range("K85").Value = 2355
$ws.Range("M85").Value = -1107
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H48").Value = 10000
$ws.Range("I48").Value = 10000
$ws.Range("K48").Value = 10000
$ws.Range("M48").Value = -9431
$ws.Range("H62").Value = 5697.143
$ws.Range("J62").Value = 5800
$ws.Range("L62").Value = 5800
$ws.Range("N62").Value = -7048
$ws.Range("H65").Value = 5697.143
$ws.Range("J65").Value = 5800
$ws.Range("L65").Value = 29000
$ws.Range("N65").Value = -35240
$ws.Range("H94").Value = 26974
$ws.Range("J94").Value = 26974
$ws.Range("L94").Value = 26974
$ws.Range("N94").Value = -28776
$ws.Range("H111").Value = 32000
$ws.Range("J111").Value = 32000
$ws.Range("L111").Value = 32000
$ws.Range("N111").Value = -40180
$ws.Range("H122").Value = 2092.6
$ws.Range("I122").Value = 1887.7778
$ws.Range("K122").Value = 5663.3334
$ws.Range("M122").Value = -3213.3334
$ws.Range("H132").Value = 1158.9796
$ws.Range("I132").Value = 880.7632
$ws.Range("K132").Value = 2642.2896
$ws.Range("M132").Value = -112.2896000000001
$ws.Range("H136").Value = 1537215.9
$ws.Range("I136").Value = 2481725.8
$ws.Range("J136").Value = 2387.375
$ws.Range("K136").Value = 7445177.399999999
$ws.Range("L136").Value = 7162.125
$ws.Range("M136").Value = -7442627.399999999
$ws.Range("N136").Value = -12262.125

Write-Output "Applied 220 cell updates across 8 sheets."